# edit.ps1 - applies the commit's changes to house_affordability.docx
# via the Word COM-interop object model.

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $r = $d.Content
    $r.Find.ClearFormatting()
    $null = $r.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $r.Find.Found) {
        throw "Find failed for: $find"
    }
}

# ---------------------------------------------------------------------
# 1) Remove the stray _GoBack bookmark after "...1997to2016)."
# ---------------------------------------------------------------------
$r = $d.Content
$null = $r.Find.Execute("_GoBack", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($r.Find.Found) {
    $bm = $d.Bookmarks("_GoBack")
    $bm.Delete()
}

# ---------------------------------------------------------------------
# 2) Merge the split hyperlink run "https://www.ons.gov.uk/people" + "p" +
#    "opulationandcommunity/.../1997to2016" into a single run.
# ---------------------------------------------------------------------
Replace-Text "https://www.ons.gov.uk/peoplepopulationandcommunity/housing/bulletins/housingaffordabilityinenglandandwales/1997to2016" "https://www.ons.gov.uk/peoplepopulationandcommunity/housing/bulletins/housingaffordabilityinenglandandwales/1997to2016"

# ---------------------------------------------------------------------
# 3) Merge " " + "Fig. 5 shows this real income per head against time. "
# ---------------------------------------------------------------------
Replace-Text " Fig. 5 shows this real income per head against time. " " Fig. 5 shows this real income per head against time. "

# ---------------------------------------------------------------------
# 4) Merge the whole "millennial's cry ... evaporate away." paragraph
#    into a single run.
# ---------------------------------------------------------------------
Replace-Text "The millennial’s cry of previous generations being lucky to buy houses when they were “cheap” does seem to hold up based on the affordability ratio alone. Older generations will refer to the punitively high interest rates in the past, but those rates were often accompanied by similar levels of inflation that helped their debts to evaporate away." "The millennial’s cry of previous generations being lucky to buy houses when they were “cheap” does seem to hold up based on the affordability ratio alone. Older generations will refer to the punitively high interest rates in the past, but those rates were often accompanied by similar levels of inflation that helped their debts to evaporate away."

# ---------------------------------------------------------------------
# 5) Merge "Houses are expensive things...one." + " And they are...because of"
# ---------------------------------------------------------------------
Replace-Text "Houses are expensive things, especially if you don’t have one. And they are getting more and more expensive all the time. Although income can be measured in many ways, the trend is clear in that wages cannot keep up with rising house prices. This makes it tough to buy your first house because of" "Houses are expensive things, especially if you don’t have one. And they are getting more and more expensive all the time. Although income can be measured in many ways, the trend is clear in that wages cannot keep up with rising house prices. This makes it tough to buy your first house because of"

# ---------------------------------------------------------------------
# 6) Merge " upfront deposit. " + "As well as...ever make." + " Houses that...forcing "
# ---------------------------------------------------------------------
Replace-Text " upfront deposit. As well as being somewhere to live, a house is the largest and least diversified investment a person will ever make. Houses that cost many multiples of their salaries are forcing " " upfront deposit. As well as being somewhere to live, a house is the largest and least diversified investment a person will ever make. Houses that cost many multiples of their salaries are forcing "

# ---------------------------------------------------------------------
# 7) "Your comments below are appreciated, both millennials and are wiser
#    forebears." -> split into 3 runs, with "are" replaced by "our".
# ---------------------------------------------------------------------
Replace-Text "Your comments below are appreciated, both millennials and are wiser forebears." "Your comments below are appreciated, both millennials and our wiser forebears."

$r = $d.Content
$null = $r.Find.Execute("our wiser forebears", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $r.Find.Found) { throw "could not find 'our wiser forebears' to split" }
$start = $r.Start
$b1 = $d.Range($start, $start)
$d.Bookmarks.Add("_tmpSplitA", $b1)
$b2 = $d.Range($start + 3, $start + 3)
$d.Bookmarks.Add("_tmpSplitB", $b2)
$d.Bookmarks("_tmpSplitA").Delete()
$d.Bookmarks("_tmpSplitB").Delete()

# ---------------------------------------------------------------------
# 8) Merge "I've contacted ONS..." + "7" + " shows a comparison..." into
#    a single run.
# ---------------------------------------------------------------------
Replace-Text "I’ve contacted ONS to understand the methodology differences between these two measures of income, so hopefully I can come back to you with an update. For now, Fig. 7 shows a comparison of the two datasets for the years they overlap. You can see how the affordability ratio differs by up to 1/3 for the two measures, however the trendlines closely match so they paint a similar picture of house prices outstripping wage growth. " "I’ve contacted ONS to understand the methodology differences between these two measures of income, so hopefully I can come back to you with an update. For now, Fig. 7 shows a comparison of the two datasets for the years they overlap. You can see how the affordability ratio differs by up to 1/3 for the two measures, however the trendlines closely match so they paint a similar picture of house prices outstripping wage growth. "

# ---------------------------------------------------------------------
# 9) "Caption: Fig. 7 Comparison of affordability ratio for median income
#    and real household's disposable income per head." -> merge "Fig. " +
#    "7" + " Comparison..." into one run, then wrap "Fig. 7 Comparison..."
#    with a new _GoBack bookmark (start right after "Caption: ").
# ---------------------------------------------------------------------
Replace-Text "Caption: Fig. 7 Comparison of affordability ratio for median income and real household’s disposable income per head." "Caption: Fig. 7 Comparison of affordability ratio for median income and real household’s disposable income per head."

$r = $d.Content
$null = $r.Find.Execute("Caption: Fig. 7 Comparison of affordability ratio for median income and real household’s disposable income per head.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $r.Find.Found) { throw "could not find caption to bookmark" }
$capStart = $r.Start + 9   # length of "Caption: "
$bmRange = $d.Range($capStart, $r.End)
$d.Bookmarks.Add("_GoBack", $bmRange)
